# Update the "Förändrad" (Changed) date column (C) from serial date 45177
# (2023-09-08) to 45178 (2023-09-09) for every data row (rows 2 through 288).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 288

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)   # Column C
    if ($cell.Value2 -eq 45177) {
        $cell.Value2 = 45178
    }
}
